$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.805.88"
$ws.Range("E2").Value = "  +4.26%  "

$ws.Range("D3").Value = "1.866.77"
$ws.Range("E3").Value = "  +2.76%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "274.23"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.64%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9979"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5273"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.64%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3386"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.22%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06819"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.00%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "19.90"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.09%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.7926"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.66%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07741"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.51%  "

$ws.Range("D13").Value = "1.855.33"
$ws.Range("E13").Value = "  +2.07%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "90.02"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.69%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.130"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.43"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.51%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008009"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.30%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.9974"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").Value = "26.840.51"
$ws.Range("E20").Value = "  +4.17%  "

$ws.Range("D21").Value = "2.107.25"
$ws.Range("E21").Value = "  +2.83%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.714"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.970"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.119"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.368"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +5.97%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "145.57"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.22%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.654"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.81%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.23"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "112.54"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.13%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.330"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.309"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.98%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.08866"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +1.39%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04923"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.90%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.163"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.64%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7275"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.873"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.223"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.24%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.335"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01846"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.25%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.5097"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.12%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.9385"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.28%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "116.34"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.00%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "6.122"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.79%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "7.981"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9970"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4415"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.50%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.1328"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.80%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.302"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.27%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "36.09"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.85%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05943"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.09%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.476"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.69%  "
